$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Trends Status" sheet (sheet1) - update values
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Trends Status")

# Row 2: Rapid Decline
$ws1.Cells.Item(2, 2).Value = 1
$ws1.Cells.Item(2, 3).Value = 0
$ws1.Cells.Item(2, 4).Value = 6.7
$ws1.Cells.Item(2, 5).Value = 0

# Row 3: Decline
$ws1.Cells.Item(3, 2).Value = 1
$ws1.Cells.Item(3, 3).Value = 8
$ws1.Cells.Item(3, 4).Value = 6.7
$ws1.Cells.Item(3, 5).Value = 21.6

# Row 4: Stable
$ws1.Cells.Item(4, 2).Value = 4
$ws1.Cells.Item(4, 4).Value = 26.7

# Row 5: Increase
$ws1.Cells.Item(5, 2).Value = 6
$ws1.Cells.Item(5, 3).Value = 4
$ws1.Cells.Item(5, 4).Value = 40
$ws1.Cells.Item(5, 5).Value = 10.8

# Row 6: Rapid Increase
$ws1.Cells.Item(6, 3).Value = 2
$ws1.Cells.Item(6, 4).Value = 20
$ws1.Cells.Item(6, 5).Value = 5.4

# Row 7: Trend Inconclusive
$ws1.Cells.Item(7, 2).Value = 9
$ws1.Cells.Item(7, 3).Value = 73

# Row 8: Insufficient Data
$ws1.Cells.Item(8, 2).Value = 435
$ws1.Cells.Item(8, 3).Value = 349

# ---------------------------------------------------------------------------
# 2. "Priority Status" sheet (sheet3) - update values
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Priority Status")
$ws3.Cells.Item(2, 2).Value = 103
$ws3.Cells.Item(3, 2).Value = 286
$ws3.Cells.Item(4, 2).Value = 554

# ---------------------------------------------------------------------------
# 3. "Species qualification" sheet (sheet4) - update values/text
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Cells.Item(2, 1).Value = "SoIB Assessment"
$ws4.Cells.Item(2, 2).Value = 459
$ws4.Cells.Item(3, 2).Value = 24
$ws4.Cells.Item(3, 3).Value = 15
$ws4.Cells.Item(4, 2).Value = 110

# ---------------------------------------------------------------------------
# 4. Rename "High Priority break-up" sheet to "Interannual update - High Pri"
#    and replace its contents with the new interannual-update data.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("High Priority break-up")
$ws5.Name = "Interannual update - High Pri"

# Overwrite the previous single "IUCN" data row with the two new rows
# ("Trend New" and an updated "IUCN" row).
$ws5.Cells.Item(2, 1).Value = "Trend New"
$ws5.Cells.Item(2, 2).Value = 67
$ws5.Cells.Item(2, 3).Value = 65
$ws5.Cells.Item(2, 4).Value = 67
$ws5.Cells.Item(2, 5).Value = 76.09999999999999

$ws5.Cells.Item(3, 1).Value = "IUCN"
$ws5.Cells.Item(3, 2).Value = 36
$ws5.Cells.Item(3, 3).Value = 35
$ws5.Cells.Item(3, 4).Value = 21
$ws5.Cells.Item(3, 5).Value = 23.9

# ---------------------------------------------------------------------------
# 5. Add a brand-new sheet "Major update - High Priority " that holds the
#    data which used to live on the "High Priority break-up" sheet. Insert
#    it right after "Interannual update - High Pri" so the sheet order
#    matches the target workbook.
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws5)
$ws6.Name = "Major update - High Priority "

$ws6.Cells.Item(1, 1).Value = "Break-up"
$ws6.Cells.Item(1, 2).Value = "High Species (no.)"
$ws6.Cells.Item(1, 3).Value = "High Species (perc.)"
$ws6.Cells.Item(1, 4).Value = "New High Species (no.)"
$ws6.Cells.Item(1, 5).Value = "New High Species (perc.)"

# Match the bold / centered header style used by the other sheets.
$hdr6 = $ws6.Range("A1:E1")
$hdr6.Font.Bold = $true
$hdr6.HorizontalAlignment = -4108

$ws6.Cells.Item(2, 1).Value = "IUCN"
$ws6.Cells.Item(2, 2).Value = 18
$ws6.Cells.Item(2, 3).Value = 100
$ws6.Cells.Item(2, 4).Value = 18
$ws6.Cells.Item(2, 5).Value = 100
